$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the example names in column B (rows 12-16) with the "real" group names
# used elsewhere in the workbook (Alice, Bob, Claire, David, Elaine), and clear
# the special font styling that was applied to those cells.
$ws.Range("B12").Value = "Alice"
$ws.Range("B13").Value = "Bob"
$ws.Range("B14").Value = "Claire"
$ws.Range("B15").Value = "David"
$ws.Range("B16").Value = "Elaine"

$ws.Range("B12:B16").Style = "Normal"

# Remove the extra row (was "Martin") entirely, shifting rows up.
$ws.Rows.Item(17).Delete()

# Update the active selection to match the new layout.
$ws.Range("D13").Select()
